$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'250.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "'5.418"
$ws.Range("D4").Style = "Normal"
$ws.Range("D6").Value = "'3.424"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.368"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8165"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9125"
$ws.Range("D9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1439"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07508"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03122"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03112"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09351"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.560"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001581"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04755"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005787"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006413"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.004998"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.001031"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'3.698"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.190"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3300"
$ws.Range("D25").Style = "Normal"
$ws.Range("E27").Value = "26AAXTokenAAB"
$ws.Range("D40").Value = "'0.04050"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006767"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1071"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002719"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007561"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005800"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.4997"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Value = "'0.01009"
$ws.Range("D50").Style = "Normal"
